$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 17.05054909069899
$ws.Range("D2").Value = 7.676382543134104
$ws.Range("E2").Value = 27.50581179381547
$ws.Range("F2").Value = 44.27141919123338
$ws.Range("G2").Value = 3.651604079620892
$ws.Range("I2").Value = 26.34611176671655
$ws.Range("L2").Value = 12.43441703384168
$ws.Range("M2").Value = 17.38399677081187
$ws.Range("B3").Value = 16.67879558153554
$ws.Range("D3").Value = 7.633035125208915
$ws.Range("E3").Value = 25.81246717813814
$ws.Range("F3").Value = 43.24329566181893
$ws.Range("G3").Value = 3.660063294624767
$ws.Range("I3").Value = 26.63002299561522
$ws.Range("L3").Value = 12.24082365066443
$ws.Range("M3").Value = 17.18722921392477
$ws.Range("B4").Value = 16.44924081419961
$ws.Range("D4").Value = 7.610667789945831
$ws.Range("E4").Value = 24.71344922119778
$ws.Range("F4").Value = 42.62499661438524
$ws.Range("G4").Value = 3.665499540223353
$ws.Range("I4").Value = 26.81199966931285
$ws.Range("L4").Value = 12.12397508631252
$ws.Range("M4").Value = 17.06896569252702
$ws.Range("B5").Value = 16.3554955109416
$ws.Range("D5").Value = 7.602614418295346
$ws.Range("E5").Value = 24.25073061485654
$ws.Range("F5").Value = 42.37667567551571
$ws.Range("G5").Value = 3.667776193795939
$ws.Range("I5").Value = 26.88809033227034
$ws.Range("L5").Value = 12.07691994136726
$ws.Range("M5").Value = 17.02146011713212
$ws.Range("B6").Value = 16.33992086773836
$ws.Range("D6").Value = 7.601341030038881
$ws.Range("E6").Value = 24.17300102663913
$ws.Range("F6").Value = 42.33567304436276
$ws.Range("G6").Value = 3.668157947437787
$ws.Range("I6").Value = 26.90084218346852
$ws.Range("L6").Value = 12.06914194859871
$ws.Range("M6").Value = 17.01361470834069
$ws.Range("B7").Value = 16.44797717187231
$ws.Range("D7").Value = 7.610554892043786
$ws.Range("E7").Value = 24.70726891817366
$ws.Range("F7").Value = 42.62163243232501
$ws.Range("G7").Value = 3.665529995050244
$ws.Range("I7").Value = 26.81301801313744
$ws.Range("L7").Value = 12.12333813962629
$ws.Range("M7").Value = 17.06832217124869
$ws.Range("B8").Value = 16.92271198167114
$ws.Range("D8").Value = 7.660551336923758
$ws.Range("E8").Value = 26.93427691759027
$ws.Range("F8").Value = 43.91445795114858
$ws.Range("G8").Value = 3.654470793615621
$ws.Range("I8").Value = 26.44242122172093
$ws.Range("L8").Value = 12.36727973674035
$ws.Range("M8").Value = 17.31564944032071
$ws.Range("B9").Value = 17.83806608874889
$ws.Range("D9").Value = 7.792483709177832
$ws.Range("E9").Value = 30.82986051423568
$ws.Range("F9").Value = 46.53593322069798
$ws.Range("G9").Value = 3.634686456783033
$ws.Range("I9").Value = 25.77599739319106
$ws.Range("L9").Value = 12.85926365890525
$ws.Range("M9").Value = 17.81884442011634
$ws.Range("B10").Value = 18.49433941609228
$ws.Range("D10").Value = 7.91015824069337
$ws.Range("E10").Value = 33.4043402720523
$ws.Range("F10").Value = 48.49256446729758
$ws.Range("G10").Value = 3.621283358729915
$ws.Range("I10").Value = 25.32255769341112
$ws.Range("L10").Value = 13.22581116544906
$ws.Range("M10").Value = 18.19681822286604
$ws.Range("B11").Value = 18.78806964313493
$ws.Range("D11").Value = 7.968164934466861
$ws.Range("E11").Value = 34.51333905119021
$ws.Range("F11").Value = 49.38492264175672
$ws.Range("G11").Value = 3.615425808032569
$ws.Range("I11").Value = 25.12400647275316
$ws.Range("L11").Value = 13.39299893861475
$ws.Range("M11").Value = 18.36996403010465
$ws.Range("B12").Value = 18.89850649441471
$ws.Range("D12").Value = 7.990768199069318
$ws.Range("E12").Value = 34.92438235336168
$ws.Range("F12").Value = 49.72281814721153
$ws.Range("G12").Value = 3.613241670679495
$ws.Range("I12").Value = 25.04992085856295
$ws.Range("L12").Value = 13.4563182062166
$ws.Range("M12").Value = 18.43565284265668
$ws.Range("B13").Value = 18.87475870909712
$ws.Range("D13").Value = 7.985871954128523
$ws.Range("E13").Value = 34.83625262440911
$ws.Range("F13").Value = 49.65005256914029
$ws.Range("G13").Value = 3.613710559469749
$ws.Range("I13").Value = 25.06582767500089
$ws.Range("L13").Value = 13.44268172764099
$ws.Range("M13").Value = 18.42150094530304
$ws.Range("B14").Value = 18.79717178062404
$ws.Range("D14").Value = 7.970011787259653
$ws.Range("E14").Value = 34.54733423384023
$ws.Range("F14").Value = 49.41272385364918
$ws.Range("G14").Value = 3.615245439207181
$ws.Range("I14").Value = 25.11788938094787
$ws.Range("L14").Value = 13.39820838251803
$ws.Range("M14").Value = 18.37536612654814
$ws.Range("B15").Value = 18.7495414910938
$ws.Range("D15").Value = 7.960379764464269
$ws.Range("E15").Value = 34.36920387407331
$ws.Range("F15").Value = 49.26734042195914
$ws.Range("G15").Value = 3.616190010006878
$ws.Range("I15").Value = 25.14992185591405
$ws.Range("L15").Value = 13.37096673453544
$ws.Range("M15").Value = 18.34712164329796
$ws.Range("B16").Value = 18.47503787727701
$ws.Range("D16").Value = 7.906456946483992
$ws.Range("E16").Value = 33.33061610962686
$ws.Range("F16").Value = 48.43426230437753
$ws.Range("G16").Value = 3.621670945192642
$ws.Range("I16").Value = 25.33568809010384
$ws.Range("L16").Value = 13.21488924207387
$ws.Range("M16").Value = 18.18552259607132
$ws.Range("B17").Value = 18.30533263549508
$ws.Range("D17").Value = 7.874519645735839
$ws.Range("E17").Value = 32.67757517893393
$ws.Range("F17").Value = 47.92352446869953
$ws.Range("G17").Value = 3.625094358674867
$ws.Range("I17").Value = 25.45162085096669
$ws.Range("L17").Value = 13.1192147277826
$ws.Range("M17").Value = 18.08665871376878
$ws.Range("B18").Value = 18.20727532055141
$ws.Range("D18").Value = 7.856571962646772
$ws.Range("E18").Value = 32.29611181259808
$ws.Range("F18").Value = 47.62999549511531
$ws.Range("G18").Value = 3.627085994622055
$ws.Range("I18").Value = 25.51902954963134
$ws.Range("L18").Value = 13.06423020090114
$ws.Range("M18").Value = 18.02991112130875
$ws.Range("B19").Value = 18.17400114426982
$ws.Range("D19").Value = 7.850567788521444
$ws.Range("E19").Value = 32.16594921445516
$ws.Range("F19").Value = 47.5306627646187
$ws.Range("G19").Value = 3.627764220261896
$ws.Range("I19").Value = 25.5419781501568
$ws.Range("L19").Value = 13.04562280211534
$ws.Range("M19").Value = 18.01071889590275
$ws.Range("B20").Value = 18.32344509755498
$ws.Range("D20").Value = 7.877875809322281
$ws.Range("E20").Value = 32.74769794570469
$ws.Range("F20").Value = 47.97787176492601
$ws.Range("G20").Value = 3.624727596942458
$ws.Range("I20").Value = 25.43920440796272
$ws.Range("L20").Value = 13.12939519715533
$ws.Range("M20").Value = 18.0971712769334
$ws.Range("B21").Value = 18.81998321338309
$ws.Range("D21").Value = 7.97465306383973
$ws.Range("E21").Value = 34.632438151364
$ws.Range("F21").Value = 49.48243629738227
$ws.Range("G21").Value = 3.61479368906486
$ws.Range("I21").Value = 25.10256777359528
$ws.Range("L21").Value = 13.41127146894975
$ws.Range("M21").Value = 18.3889141244427
$ws.Range("B22").Value = 19.13983924455314
$ws.Range("D22").Value = 8.041613523530863
$ws.Range("E22").Value = 35.81230503480973
$ws.Range("F22").Value = 50.46548655778771
$ws.Range("G22").Value = 3.608499205343521
$ws.Range("I22").Value = 24.88897215171618
$ws.Range("L22").Value = 13.59552163629018
$ws.Range("M22").Value = 18.58027455436396
$ws.Range("B23").Value = 18.96958382821223
$ws.Range("D23").Value = 8.005538483917126
$ws.Range("E23").Value = 35.18732803206861
$ws.Range("F23").Value = 49.94094837564175
$ws.Range("G23").Value = 3.611840736749996
$ws.Range("I23").Value = 25.00238800390025
$ws.Range("L23").Value = 13.49719895001573
$ws.Range("M23").Value = 18.4780953704834
$ws.Range("B24").Value = 18.3152579767031
$ws.Range("D24").Value = 7.876357199364234
$ws.Range("E24").Value = 32.71601420015622
$ws.Range("F24").Value = 47.95330101098138
$ws.Range("G24").Value = 3.624893336782864
$ws.Range("I24").Value = 25.44481552201396
$ws.Range("L24").Value = 13.12479254039128
$ws.Range("M24").Value = 18.09241826115674
$ws.Range("B25").Value = 17.59282653020763
$ws.Range("D25").Value = 7.753146386322217
$ws.Range("E25").Value = 29.82655572448443
$ws.Range("F25").Value = 45.81982399594403
$ws.Range("G25").Value = 3.639837763827555
$ws.Range("I25").Value = 25.94988635297622
$ws.Range("L25").Value = 12.7250447923863
$ws.Range("M25").Value = 17.68105855489277
